# Apply the edits described by the diff:
#  - remove the last "skip" row (row 21, "Пропустить") from the list
#  - add a thin border around the whole table (A1:F20)
#  - give the header row (A1:F1) a light green (Accent6, lighter 60%) fill
#  - add conditional formatting on D2:F20 that highlights values > 0
#    with the standard "Green Fill with Dark Green Text" look
#  - select cell A21 (first empty row after the table) as the active cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra "Пропустить" (skip) row at the bottom of the table.
$ws.Rows.Item(21).Delete() | Out-Null

# Add a thin box border around every cell of the table (header + data).
$table = $ws.Range("A1:F20")
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2

# Highlight the header row with a light green fill (theme Accent6,
# lighter 60% -> RGB C5E0B4, passed in BGR order for the COM Color prop).
$headerRow = $ws.Range("A1:F1")
$headerRow.Interior.Color = 0xB4E0C5

# Conditional formatting: values greater than 0 in D2:F20 get the
# standard "green fill / dark green text" highlight.
$dataRange = $ws.Range("D2:F20")
$condition = $dataRange.FormatConditions.Add(1, 5, "0")
$condition.Font.Color = 0x006100
$condition.Interior.Color = 0xCEEFC6

# Select A21 (the now-empty row right after the table) to match the
# saved selection state of the workbook.
$ws.Range("A21").Select() | Out-Null
